$d = $word.ActiveDocument

# --- Paragraph 1 (the "**ID__...**" marker paragraph) ---
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat

# Add a paragraph border (all four sides) with 5pt spacing and no visible line.
$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Bump the left indent from 6pt (120 twips) to 11.25pt (225 twips).
$pf.LeftIndent = 11.25

# Replace the marker text and swallow the trailing space run that followed it,
# turning "**ID__AFFARS_5350_topic_7__ID** " into "**ID__AFFARS_SUBPART_5350_103__ID**".
$find = $p1.Range.Find
[void]$find.Execute("**ID__AFFARS_5350_topic_7__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5350_103__ID**", 2)
